$wb = $excel.ActiveWorkbook

# --- Add the new "Trait codings" sheet as the last tab ---
$newSheet = $wb.Worksheets.Add()
$newSheet.Name = "Trait codings"
$newSheet.Move($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$ws3 = $wb.Worksheets.Item("Trait codings")

# --- Populate the human readable trait codings table ---
$ws3.Range("A1").Value = "Trait name"
$ws3.Range("B1").Value = "Long name"
$ws3.Range("C1").Value = 0
$ws3.Range("D1").Value = 1
$ws3.Range("E1").Value = 2
$ws3.Range("A2").Value = "Size"
$ws3.Range("B2").Value = "Valve Size"
$ws3.Range("C2").Value = "Small (< 0.50mm x 0.25mm)"
$ws3.Range("D2").Value = "Medium (0.50mm-1.00mm x 0.25mm – 0.50mm)"
$ws3.Range("E2").Value = "Large (> 1.00mm x 0.50mm)"
$ws3.Range("A3").Value = "Calcification"
$ws3.Range("B3").Value = "Valve Calcification"
$ws3.Range("C3").Value = "Weak"
$ws3.Range("D3").Value = "Average"
$ws3.Range("E3").Value = "Thin"
$ws3.Range("A4").Value = "Shape"
$ws3.Range("B4").Value = "Valve Shape (dorsal view: anterior-posterior)"
$ws3.Range("C4").Value = "Ovate or Subovate"
$ws3.Range("D4").Value = "Rectangular or Subrectangular"
$ws3.Range("A5").Value = "Carapace Texture"
$ws3.Range("B5").Value = "Presence of Surface Reticulation"
$ws3.Range("C5").Value = "Absent"
$ws3.Range("D5").Value = "Present"
$ws3.Range("A6").Value = "Ventral Margin"
$ws3.Range("B6").Value = "Ventral Margin Shape"
$ws3.Range("C6").Value = "Convex"
$ws3.Range("D6").Value = "Straight"
$ws3.Range("E6").Value = "Concave"
$ws3.Range("A7").Value = "Dorsal Margin"
$ws3.Range("B7").Value = "Dorsal Margin Shape"
$ws3.Range("C7").Value = "Convex"
$ws3.Range("D7").Value = "Straight"
$ws3.Range("E7").Value = "Concave"
$ws3.Range("A8").Value = "Posterior Margin"
$ws3.Range("B8").Value = "Posterior Margin Shape"
$ws3.Range("C8").Value = "Convex"
$ws3.Range("D8").Value = "Straight"
$ws3.Range("E8").Value = "Concave"
$ws3.Range("A9").Value = "Anterior Margin"
$ws3.Range("B9").Value = "Anterior Margin Shape"
$ws3.Range("C9").Value = "Convex"
$ws3.Range("D9").Value = "Straight"
$ws3.Range("E9").Value = "Concave"
$ws3.Range("A10").Value = "Left Overlap"
$ws3.Range("B10").Value = "Right/Left Valve Size Ratio"
$ws3.Range("C10").Value = "Left Valve Larger"
$ws3.Range("D10").Value = "Equally Sized Valves"
$ws3.Range("E10").Value = "Right Valve Larger"
$ws3.Range("A11").Value = "Spines"
$ws3.Range("B11").Value = "Presence of Spines"
$ws3.Range("C11").Value = "Absent"
$ws3.Range("D11").Value = "Present"
$ws3.Range("A12").Value = "Carapace Pits"
$ws3.Range("B12").Value = "Presence of Carapace Pits"
$ws3.Range("C12").Value = "Absent"
$ws3.Range("D12").Value = "Present"
$ws3.Range("A13").Value = "Opaque Areas"
$ws3.Range("B13").Value = "Presence of Opaque Patches on Valve"
$ws3.Range("C13").Value = "Absent"
$ws3.Range("D13").Value = "Present"
$ws3.Range("A14").Value = "Denticulations"
$ws3.Range("B14").Value = "Presence of Denticulations"
$ws3.Range("C14").Value = "Absent"
$ws3.Range("D14").Value = "Present"
$ws3.Range("A15").Value = "ala"
$ws3.Range("B15").Value = "Presence of Alae"
$ws3.Range("C15").Value = "Absent"
$ws3.Range("D15").Value = "Present"
$ws3.Range("A16").Value = "nodes"
$ws3.Range("B16").Value = "Presence of Nodes"
$ws3.Range("C16").Value = "Absent"
$ws3.Range("D16").Value = "Present"
$ws3.Range("A17").Value = "caudal process"
$ws3.Range("B17").Value = "Presence of Caudal Process"
$ws3.Range("C17").Value = "Absent"
$ws3.Range("D17").Value = "Present"
$ws3.Range("A18").Value = "sulcus"
$ws3.Range("B18").Value = "Presence of Sulcus"
$ws3.Range("C18").Value = "Absent"
$ws3.Range("D18").Value = "Present"
$ws3.Range("A19").Value = "eye tubercule"
$ws3.Range("B19").Value = "Presence of Eye Tubercules"
$ws3.Range("C19").Value = "Absent"
$ws3.Range("D19").Value = "Present"

# --- Traits sheet: scroll/selection reset (was scrolled to A110 / A3 selected) ---
$ws1 = $wb.Worksheets.Item("Traits")
$ws1.Rows.Item(35).RowHeight = 13.8
$ws1.Range("A2").Select()

# --- Activate the new sheet last so it becomes the selected tab ---
$ws3.Activate()
$ws3.Range("A12").Select()
